# Remove Defect ID (DE...) values from the RTM sheet's "Comments Failed US-
# Defect ID Not Testable Comment" column (J), per the Release Agent's request.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTM")

$rows = @(3, 4, 5, 21, 23, 28, 29, 30, 31, 32, 33, 34, 35, 36, 40)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 10).Value = ""
}

# Leave the sheet's selection where the editor was last working (last
# defect-id cell that got cleared out) rather than the original G6.
$ws.Activate() | Out-Null
$ws.Range("J40").Select() | Out-Null
